$wb = $excel.ActiveWorkbook

# --- Rename the "Include from SNOMED CT" sheets ---
$wb.Worksheets.Item("Include from SNOMED CT").Name = "Include #0"
$wb.Worksheets.Item("Include from SNOMED CT 2").Name = "Include #1"
$wb.Worksheets.Item("Include from SNOMED CT 3").Name = "Include #2"

# --- Update the Metadata sheet ---
$ws = $wb.Worksheets.Item("Metadata")

# Row 8: Date value updated
$ws.Range("B8").Value = "2024-09-17T19:55:11+00:00"

# Row 11: new "Jurisdiction" property (empty value)
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Row 12: Description
$ws.Range("A12").Value = "Description"
$ws.Range("B12").Value = "Result values for N category. This value set contains SNOMED-CT equivalents of UICC codes for the N category, according to TNM staging rules."

# Row 13: Purpose (value stays empty)
$ws.Range("A13").Value = "Purpose"
$ws.Range("B13").Value = ""

# Row 14: Copyright
$ws.Range("A14").Value = "Copyright"
$ws.Range("B14").Value = "This value set includes content from SNOMED CT, which is copyright © 2002+ International Health Terminology Standards Development Organisation (IHTSDO), and distributed by agreement between IHTSDO and HL7. Implementer use of SNOMED CT is not covered by this agreement"

# Row 15: new "Immutable" row - copy formatting from row 13 first so the
# new cells reuse the existing data-row style instead of creating a new one
$ws.Range("A13").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("B13").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("A15").Value = "Immutable"
$ws.Range("B15").Value = "BooleanType[null]"

Write-Host "Edit applied"
